$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 7830
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = 1202
$ws.Range("G2").Value = 1127
$ws.Range("H2").Value = 850
$ws.Range("I2").Value = 449
$ws.Range("J2").Value = 401
$ws.Range("K2").Value = 53385
$ws.Range("L2").Value = 40884
$ws.Range("M2").Value = 12501
$ws.Range("N2").Value = 7601
$ws.Range("O2").Value = 4900
$ws.Range("P2").Value = 224
$ws.Range("Q2").Value = -2049
$ws.Range("R2").Value = -591
$ws.Range("S2").Value = 2633
$ws.Range("T2").Value = 102
$ws.Range("U2").Value = -2151
$ws.Range("V2").Value = 9257
$ws.Range("W2").Value = 15.36
$ws.Range("X2").Value = 10.85
$ws.Range("Y2").Value = 6.06
$ws.Range("Z2").Value = 1.67
$ws.Range("AA2").Value = 327.06
$ws.Range("AB2").Value = 3377.49
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 11.95
$ws.Range("AE2").Value = 17567
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 130
$ws.Range("AH2").Value = 1.09
$ws.Range("AI2").Value = 12.54
$ws.Range("AJ2").Value = 44866617

# Row 3
$ws.Range("D3").Value = 10385
$ws.Range("E3").Value = 2605
$ws.Range("F3").Value = 2605
$ws.Range("G3").Value = 2831
$ws.Range("H3").Value = 2026
$ws.Range("I3").Value = 1031
$ws.Range("J3").Value = 995
$ws.Range("K3").Value = 68347
$ws.Range("L3").Value = 53871
$ws.Range("M3").Value = 14476
$ws.Range("N3").Value = 8507
$ws.Range("O3").Value = 5969
$ws.Range("P3").Value = 224
$ws.Range("Q3").Value = -1614
$ws.Range("R3").Value = -1394
$ws.Range("S3").Value = 3257
$ws.Range("T3").Value = 111
$ws.Range("U3").Value = -1725
$ws.Range("V3").Value = 14400
$ws.Range("W3").Value = 25.08
$ws.Range("X3").Value = 19.51
$ws.Range("Y3").Value = 12.8
$ws.Range("Z3").Value = 3.33
$ws.Range("AA3").Value = 372.15
$ws.Range("AB3").Value = 3781.8
$ws.Range("AC3").Value = 2297
$ws.Range("AD3").Value = 9.49
$ws.Range("AE3").Value = 19661
$ws.Range("AF3").Value = 1.11
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 0.6899999999999999
$ws.Range("AI3").Value = 6.3
$ws.Range("AJ3").Value = 44866617

# Row 4
$ws.Range("D4").Value = 11534
$ws.Range("E4").Value = 2534
$ws.Range("F4").Value = 2534
$ws.Range("G4").Value = 2608
$ws.Range("H4").Value = 1920
$ws.Range("I4").Value = 975
$ws.Range("J4").Value = 945
$ws.Range("K4").Value = 95446
$ws.Range("L4").Value = 79370
$ws.Range("M4").Value = 16076
$ws.Range("N4").Value = 9310
$ws.Range("O4").Value = 6766
$ws.Range("P4").Value = 224
$ws.Range("Q4").Value = -3574
$ws.Range("R4").Value = -6318
$ws.Range("S4").Value = 9572
$ws.Range("T4").Value = 123
$ws.Range("U4").Value = -3698
$ws.Range("V4").Value = 30714
$ws.Range("W4").Value = 21.97
$ws.Range("X4").Value = 16.64
$ws.Range("Y4").Value = 10.95
$ws.Range("Z4").Value = 2.34
$ws.Range("AA4").Value = 493.73
$ws.Range("AB4").Value = 4198.59
$ws.Range("AC4").Value = 2173
$ws.Range("AD4").Value = 8.279999999999999
$ws.Range("AE4").Value = 21517
$ws.Range("AF4").Value = 0.84
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 1.39
$ws.Range("AI4").Value = 11.09
$ws.Range("AJ4").Value = 44866617

# Row 5
$ws.Range("D5").Value = 13728
$ws.Range("E5").Value = 3302
$ws.Range("F5").Value = 3302
$ws.Range("G5").Value = 3438
$ws.Range("H5").Value = 2431
$ws.Range("I5").Value = 1146
$ws.Range("J5").Value = 1286
$ws.Range("K5").Value = 125640
$ws.Range("L5").Value = 106531
$ws.Range("M5").Value = 19110
$ws.Range("N5").Value = 10570
$ws.Range("O5").Value = 8539
$ws.Range("P5").Value = 224
$ws.Range("Q5").Value = -8514
$ws.Range("R5").Value = -2734
$ws.Range("S5").Value = 11567
$ws.Range("T5").Value = 109
$ws.Range("U5").Value = -8622
$ws.Range("V5").Value = 45668
$ws.Range("W5").Value = 24.05
$ws.Range("X5").Value = 17.71
$ws.Range("Y5").Value = 11.53
$ws.Range("Z5").Value = 2.2
$ws.Range("AA5").Value = 557.46
$ws.Range("AB5").Value = 4666.52
$ws.Range("AC5").Value = 2554
$ws.Range("AD5").Value = 7.54
$ws.Range("AE5").Value = 24431
$ws.Range("AF5").Value = 0.79
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 1.56
$ws.Range("AI5").Value = 11.33
$ws.Range("AJ5").Value = 44866617

# Row 6
$ws.Range("D6").Value = 23546
$ws.Range("E6").Value = 3239
$ws.Range("F6").Value = 3239
$ws.Range("G6").Value = 3244
$ws.Range("H6").Value = 2254
$ws.Range("I6").Value = 1104
$ws.Range("K6").Value = 189210
$ws.Range("L6").Value = 164804
$ws.Range("M6").Value = 24406
$ws.Range("N6").Value = 11999
$ws.Range("P6").Value = 224
$ws.Range("Q6").Value = -14414
$ws.Range("R6").Value = -7839
$ws.Range("S6").Value = 25286
$ws.Range("T6").Value = 272
$ws.Range("U6").Value = -14686
$ws.Range("V6").Value = 64049
$ws.Range("W6").Value = 13.76
$ws.Range("X6").Value = 9.57
$ws.Range("Y6").Value = 9.779999999999999
$ws.Range("Z6").Value = 1.43
$ws.Range("AA6").Value = 675.27
$ws.Range("AB6").Value = 5216.01
$ws.Range("AC6").Value = 2460
$ws.Range("AD6").Value = 7.66
$ws.Range("AE6").Value = 27732
$ws.Range("AF6").Value = 0.68
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 1.86
$ws.Range("AI6").Value = 13.72
$ws.Range("AJ6").Value = 44866617

# Clear rows 7-9 data cells (D:AJ), keep A/B/C headers intact
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
